$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 985
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 985
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 985
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -1553
# Row 19
$ws.Range("H19").Value = 1411.1111
$ws.Range("J19").Value = 1462.625
$ws.Range("L19").Value = 1462.625
$ws.Range("N19").Value = -1812.625
# Row 40
$ws.Range("H40").Value = 3810.3044
$ws.Range("J40").Value = 4183.278
$ws.Range("L40").Value = 4183.278
$ws.Range("N40").Value = -4533.278
# Row 69
$ws.Range("H69").Value = 21701.75
$ws.Range("I69").Value = 5896
$ws.Range("J69").Value = 37507.5
$ws.Range("K69").Value = 17688
$ws.Range("L69").Value = 112522.5
$ws.Range("M69").Value = -16814
$ws.Range("N69").Value = -114270.5
# Row 70
$ws.Range("H70").Value = 1659.6
$ws.Range("I70").Value = 1659.6
$ws.Range("K70").Value = 4978.799999999999
$ws.Range("M70").Value = -4708.799999999999
# Row 72
$ws.Range("H72").Value = 21701.75
$ws.Range("I72").Value = 5896
$ws.Range("J72").Value = 37507.5
$ws.Range("K72").Value = 53064
$ws.Range("L72").Value = 337567.5
$ws.Range("M72").Value = -48696
$ws.Range("N72").Value = -346303.5
# Row 73
$ws.Range("H73").Value = 1659.6
$ws.Range("I73").Value = 1659.6
$ws.Range("K73").Value = 4978.799999999999
$ws.Range("M73").Value = -4042.799999999999
# Row 101
$ws.Range("H101").Value = 3518.2
$ws.Range("I101").Value = 3518.2
$ws.Range("K101").Value = 10554.6
$ws.Range("M101").Value = -8932.599999999999
# Row 138
$ws.Range("H138").Value = 2061.739
$ws.Range("J138").Value = 2593.2334
$ws.Range("L138").Value = 7779.7002
$ws.Range("N138").Value = -18059.7002
# Row 141
$ws.Range("H141").Value = 3069.6667
$ws.Range("I141").Value = 3069.6667
$ws.Range("K141").Value = 9209.000100000001
$ws.Range("M141").Value = -4029.000100000001

$ws = $wb.Worksheets.Item("ARM")
# Row 81
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").ClearContents()
# Row 84
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").ClearContents()
# Row 124
$ws.Range("H124").Value = 80167.8
$ws.Range("J124").Value = 80167.8
$ws.Range("L124").Value = 80167.8
$ws.Range("N124").Value = -89987.8
# Row 125
$ws.Range("H125").Value = 47996.668
$ws.Range("J125").Value = 47996.668
$ws.Range("L125").Value = 47996.668
$ws.Range("N125").Value = -57836.668
# Row 130
$ws.Range("H130").Value = 78486
$ws.Range("J130").Value = 78486
$ws.Range("L130").Value = 78486
$ws.Range("N130").Value = -88526
# Row 135
$ws.Range("H135").Value = 62999.5
$ws.Range("J135").Value = 62999.5
$ws.Range("L135").Value = 62999.5
$ws.Range("N135").Value = -73139.5

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1496.9
$ws.Range("I94").Value = 1839.2858
$ws.Range("J94").Value = 698
$ws.Range("K94").Value = 1839.2858
$ws.Range("L94").Value = 698
$ws.Range("M94").Value = -1388.2858
$ws.Range("N94").Value = -1600
# Row 105
$ws.Range("H105").Value = 2551.9412
$ws.Range("I105").Value = 2034.9
$ws.Range("K105").Value = 2034.9
$ws.Range("M105").Value = -287.9000000000001
# Row 132
$ws.Range("H132").Value = 83500
$ws.Range("J132").Value = 83500
$ws.Range("L132").Value = 83500
$ws.Range("N132").Value = -93620

$ws = $wb.Worksheets.Item("CRP")
# Row 134
$ws.Range("H134").Value = 348734.1
$ws.Range("I134").Value = 456496.4
$ws.Range("J134").Value = 10052.571
$ws.Range("K134").Value = 1369489.2
$ws.Range("L134").Value = 30157.713
$ws.Range("M134").Value = -1366954.2
$ws.Range("N134").Value = -35227.713

$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 531.125
$ws.Range("I8").Value = 531.125
$ws.Range("K8").Value = 1593.375
$ws.Range("M8").Value = -1454.375
# Row 23
$ws.Range("H23").Value = 67.5
$ws.Range("J23").Value = 67.5
$ws.Range("L23").Value = 202.5
$ws.Range("N23").Value = -672.5
# Row 46
$ws.Range("H46").Value = 1094.4
$ws.Range("J46").Value = 1559.2858
$ws.Range("L46").Value = 4677.857400000001
$ws.Range("N46").Value = -4859.857400000001
# Row 122
$ws.Range("H122").Value = 2120.75
$ws.Range("J122").Value = 2427.6667
$ws.Range("L122").Value = 21849.0003
$ws.Range("N122").Value = -26749.0003

$ws = $wb.Worksheets.Item("GSM")
# Row 22
$ws.Range("H22").Value = 5007
$ws.Range("I22").Value = 5007
$ws.Range("K22").Value = 5007
$ws.Range("M22").Value = -4478
# Row 35
$ws.Range("H35").Value = 90560
$ws.Range("I35").Value = 109000
$ws.Range("K35").Value = 109000
$ws.Range("M35").Value = -108702
# Row 125
$ws.Range("H125").Value = 87247.25
$ws.Range("J125").Value = 87247.25
$ws.Range("L125").Value = 87247.25
$ws.Range("N125").Value = -92167.25

$ws = $wb.Worksheets.Item("LTW")
# Row 81
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").ClearContents()
# Row 84
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").ClearContents()
# Row 127
$ws.Range("H127").Value = 86954
$ws.Range("J127").Value = 86954
$ws.Range("L127").Value = 86954
$ws.Range("N127").Value = -96874
# Row 136
$ws.Range("H136").Value = 47513.816
$ws.Range("I136").Value = 5094.7
$ws.Range("K136").Value = 15284.1
$ws.Range("M136").Value = -12734.1

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 70333.664
$ws.Range("I81").Value = 999
$ws.Range("K81").Value = 1998
$ws.Range("M81").Value = -937
# Row 84
$ws.Range("H84").Value = 70333.664
$ws.Range("I84").Value = 999
$ws.Range("K84").Value = 9990
$ws.Range("M84").Value = -4686
# Row 135
$ws.Range("H135").Value = 63999.5
$ws.Range("J135").Value = 63999.5
$ws.Range("L135").Value = 63999.5
$ws.Range("N135").Value = -74139.5
